$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 238, shifting existing rows 238..365 down to 239..366.
$ws.Rows(238).Insert()

# Populate the newly inserted row 238 with the new weekly data point.
$ws.Range("A238").Value = 10
$ws.Range("B238").Value = "Vega Modelo de Temuco"
$ws.Range("C238").Value = "La Araucanía"
$ws.Range("D238").Value = 44813
$ws.Range("E238").Value = 9
$ws.Range("F238").Value = 100112017
$ws.Range("G238").Value = "Apio"
$ws.Range("H238").Value = "Americana (o)"
$ws.Range("I238").Value = "Primera"
$ws.Range("J238").Value = 50
$ws.Range("K238").Value = 12000
$ws.Range("L238").Value = 12000
$ws.Range("M238").Value = 12000
$ws.Range("N238").Value = "$/docena de matas"
$ws.Range("O238").Value = "Provincia del Elquí"
$ws.Range("P238").Value = 2000
$ws.Range("Q238").Value = 6
$ws.Range("R238").Value = "Hortaliza"
